# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
#
# The sheet currently has data in columns A:AC across rows 1 (header) and
# 2:55 (one row per player). We extend it with three new trailing columns:
#   AD -> "Wins"
#   AE -> "Losses"
#   AF -> "Ties"
# using the same header style as the existing header cells (bold, bordered,
# centered) and fill every data row with the team's season record
# (74 wins, 88 losses, 0 ties).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the three
# new header cells so they pick up the same style (bold font, border,
# centered alignment) instead of the default style.
$headerFormat = $ws.Range("AC1")
$newHeaders = $ws.Range("AD1:AF1")
$headerFormat.Copy($newHeaders)

# Now set the actual header text (after the format copy, so these values are
# not clobbered by the copied-over content of AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every player row (2 through 55).
$lastRow = 55
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 74   # column AD - Wins
    $ws.Cells.Item($r, 31).Value = 88   # column AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # column AF - Ties
}
